# "update DORA for better rendering" / "remove ref_id for articles"
#
# 1) Clear the Ref_id (column C) values for the "article" rows on the
#    dora sheet, keeping the cell (and its style) in place.
# 2) Update the saved view state: the "dora" sheet becomes the active /
#    selected tab (scrolled back to the top, cursor on C297), while the
#    "library_content" sheet is no longer the selected tab and is left
#    scrolled down with the cursor on B11.

$wb = $excel.ActiveWorkbook

$wsLibrary = $wb.Worksheets.Item("library_content")
$wsDora    = $wb.Worksheets.Item("dora")

# --- 1) Remove Ref_id values for the listed article rows -------------------
$rows = @(4, 23, 42, 48, 56, 73, 79, 99, 115, 128, 132, 146, 156, 165, 182, 185, 192, 196, 211, 225, 259, 269, 297)

foreach ($r in $rows) {
    $wsDora.Range("C" + $r).ClearContents()
}

# --- 2) View state -----------------------------------------------------
# Leave "library_content" scrolled to row 13 with B11 selected (it will
# no longer be the active/selected tab once "dora" is activated below).
$wsLibrary.Activate()
$wsLibrary.Range("B11").Select()
$excel.ActiveWindow.ScrollRow = 13

# Make "dora" the active tab, scrolled back to the top with C297 selected.
$wsDora.Activate()
$wsDora.Range("C297").Select()
$excel.ActiveWindow.ScrollRow = 1
